# Auto-generated Excel COM-interop script
# Updates market-price / profit columns (H:N) across multiple sheets
# to reflect a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 11655
$ws.Range("I7").Value = 2482.5
$ws.Range("J7").Value = 30000
$ws.Range("K7").Value = 2482.5
$ws.Range("L7").Value = 30000
$ws.Range("M7").Value = -2370.5
$ws.Range("N7").Value = -30224

$ws.Range("H9").Value = 892.63635
$ws.Range("I9").Value = 892.63635
$ws.Range("K9").Value = 892.63635
$ws.Range("M9").Value = -723.63635

$ws.Range("H14").Value = 11655
$ws.Range("I14").Value = 2482.5
$ws.Range("J14").Value = 30000
$ws.Range("K14").Value = 2482.5
$ws.Range("L14").Value = 30000
$ws.Range("M14").Value = -2291.5
$ws.Range("N14").Value = -30382

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H132").Value = 1328.15
$ws.Range("I132").Value = 1257.2941
$ws.Range("J132").Value = 1729.6666
$ws.Range("K132").Value = 3771.8823
$ws.Range("L132").Value = 5188.9998
$ws.Range("M132").Value = -1241.8823
$ws.Range("N132").Value = -10248.9998

$ws.Range("H137").Value = 4387.923
$ws.Range("I137").Value = 1551.3334
$ws.Range("J137").Value = 5238.9
$ws.Range("K137").Value = 4654.0002
$ws.Range("L137").Value = 15716.7
$ws.Range("M137").Value = -2104.0002
$ws.Range("N137").Value = -20816.7

$ws.Range("H138").Value = 4726.4443
$ws.Range("J138").Value = 5112.3
$ws.Range("L138").Value = 15336.9
$ws.Range("N138").Value = -25616.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 876.0909
$ws.Range("I2").Value = 861.8461
$ws.Range("K2").Value = 861.8461
$ws.Range("M2").Value = -748.8461

$ws.Range("H32").Value = 4187.7715
$ws.Range("I32").Value = 3839.3667
$ws.Range("K32").Value = 3839.3667
$ws.Range("M32").Value = -3552.3667

$ws.Range("H61").Value = 2998.75
$ws.Range("I61").Value = 2998
$ws.Range("K61").Value = 2998
$ws.Range("M61").Value = -2786

$ws.Range("H74").Value = 9522023
$ws.Range("I74").Value = 12496093
$ws.Range("K74").Value = 12496093
$ws.Range("M74").Value = -12495219

$ws.Range("H77").Value = 9522023
$ws.Range("I77").Value = 12496093
$ws.Range("K77").Value = 62480465
$ws.Range("M77").Value = -62476097

$ws.Range("H116").Value = 876.0909
$ws.Range("I116").Value = 861.8461
$ws.Range("K116").Value = 861.8461
$ws.Range("M116").Value = 1432.1539

$ws.Range("H122").Value = 1422.2727
$ws.Range("I122").Value = 1422.2727
$ws.Range("K122").Value = 4266.8181
$ws.Range("M122").Value = -1816.8181

$ws.Range("H136").Value = 2998.75
$ws.Range("I136").Value = 2998
$ws.Range("K136").Value = 8994
$ws.Range("M136").Value = -6444

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 876.0909
$ws.Range("I3").Value = 861.8461
$ws.Range("K3").Value = 861.8461
$ws.Range("M3").Value = -747.8461

$ws.Range("H134").Value = 1604.125
$ws.Range("I134").Value = 1444.4
$ws.Range("K134").Value = 4333.200000000001
$ws.Range("M134").Value = -1798.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws.Range("H31").Value = 2564.5557
$ws.Range("I31").Value = 2843.8333
$ws.Range("J31").Value = 2006
$ws.Range("K31").Value = 2843.8333
$ws.Range("L31").Value = 2006
$ws.Range("M31").Value = -2548.8333
$ws.Range("N31").Value = -2596

$ws.Range("H34").Value = 2564.5557
$ws.Range("I34").Value = 2843.8333
$ws.Range("J34").Value = 2006
$ws.Range("K34").Value = 2843.8333
$ws.Range("L34").Value = 2006
$ws.Range("M34").Value = -2641.8333
$ws.Range("N34").Value = -2410

$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()

$ws.Range("H132").Value = 1793.3529
$ws.Range("I132").Value = 1040.9166
$ws.Range("K132").Value = 3122.7498
$ws.Range("M132").Value = -592.7498000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 134.5
$ws.Range("J38").Value = 162.66667
$ws.Range("L38").Value = 488.00001
$ws.Range("N38").Value = -1182.00001

$ws.Range("H125").Value = 7500
$ws.Range("I125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("M125").ClearContents()

$ws.Range("H138").Value = 5033.6665
$ws.Range("I138").Value = 4240.4
$ws.Range("K138").Value = 12721.2
$ws.Range("M138").Value = -7581.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1488.5
$ws.Range("J122").Value = 2203.5
$ws.Range("L122").Value = 6610.5
$ws.Range("N122").Value = -11510.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 7563
$ws.Range("I20").Value = 120
$ws.Range("K20").Value = 120
$ws.Range("M20").Value = 106

$ws.Range("H68").Value = 3000
$ws.Range("I68").Value = 2000
$ws.Range("J68").Value = 3666.6667
$ws.Range("K68").Value = 2000
$ws.Range("L68").Value = 3666.6667
$ws.Range("M68").Value = -1251
$ws.Range("N68").Value = -5164.6667

$ws.Range("H71").Value = 3000
$ws.Range("I71").Value = 2000
$ws.Range("J71").Value = 3666.6667
$ws.Range("K71").Value = 10000
$ws.Range("L71").Value = 18333.3335
$ws.Range("M71").Value = -6256
$ws.Range("N71").Value = -25821.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 60105
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

$ws.Range("H73").Value = 60105
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

$ws.Range("H107").Value = 394.42856
$ws.Range("I107").Value = 468
$ws.Range("J107").Value = 210.5
$ws.Range("K107").Value = 1404
$ws.Range("L107").Value = 631.5
$ws.Range("M107").Value = 516
$ws.Range("N107").Value = -4471.5

$ws.Range("H136").Value = 2027.3334
$ws.Range("I136").Value = 1548.8334
$ws.Range("K136").Value = 4646.5002
$ws.Range("M136").Value = -2096.5002
